$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.425.11"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.763.47"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("E4").Value = "  +0.92%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.23"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("E6").Value = "  +0.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3848"
$ws.Range("E7").Value = "  +1.67%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3409"
$ws.Range("E8").Value = "  -1.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.07"
$ws.Range("E9").Value = "  -2.73%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.138"
$ws.Range("E10").Value = "  -5.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07409"
$ws.Range("E11").Value = "  -1.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.007"
$ws.Range("E12").Value = "  +0.61%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.95"
$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.341"
$ws.Range("E14").Value = "  -2.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.767.64"
$ws.Range("E15").Value = "  -1.19%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.047"
$ws.Range("E16").Value = "  -0.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001072"
$ws.Range("E17").Value = "  -2.65%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06683"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.06"
$ws.Range("E19").Value = "  -2.36%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"
$ws.Range("E20").Value = "  +0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.33"
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.430"
$ws.Range("E22").Value = "  -3.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.475.37"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.15"
$ws.Range("E24").Value = "  -2.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.390"
$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.71"
$ws.Range("E26").Value = "  -3.33%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.431"
$ws.Range("E27").Value = "  -5.15%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.401"
$ws.Range("E28").Value = "  -7.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.55"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.49"
$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.966.22"
$ws.Range("E31").Value = "  -1.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.086"
$ws.Range("E32").Value = "  -0.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.964"
$ws.Range("E33").Value = "  -1.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08745"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.61"
$ws.Range("E35").Value = "  -5.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02419"
$ws.Range("E36").Value = "  +3.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.368"
$ws.Range("E37").Value = "  -2.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6752"
$ws.Range("E38").Value = "  -3.13%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06319"
$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2187"
$ws.Range("E40").Value = "  -1.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.538"
$ws.Range("E41").Value = "  -7.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.244"
$ws.Range("E42").Value = "  -0.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.381"
$ws.Range("E43").Value = "  -5.37%  "

$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("E45").Value = "  +0.66%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6223"
$ws.Range("E46").Value = "  -4.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.840"
$ws.Range("E47").Value = "  -0.25%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.74"
$ws.Range("E48").Value = "  +1.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.087"
$ws.Range("E49").Value = "  -3.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07359"
$ws.Range("E50").Value = "  +2.89%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.233"
$ws.Range("E51").Value = "  +1.71%  "
